# Updated cryptos list on Wed Apr 26 10:35:59 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

Set-TextCell 'D2' '28.851.53'
Set-TextCell 'E2' '  +5.07%  '
Set-TextCell 'D3' '1.904.72'
Set-TextCell 'E3' '  +4.30%  '
Set-TextCell 'D4' '1.000'
Set-TextCell 'E4' '  -0.63%  '
Set-TextCell 'D5' '338.28'
Set-TextCell 'E5' '  +1.79%  '
Set-TextCell 'D6' '0.9997'
Set-TextCell 'E6' '  -0.59%  '
Set-TextCell 'D7' '0.4711'
Set-TextCell 'E7' '  +3.02%  '
Set-TextCell 'E8' '  +6.31%  '
Set-TextCell 'D9' '47.91'
Set-TextCell 'E9' '  +3.16%  '
Set-TextCell 'D10' '0.08153'
Set-TextCell 'E10' '  +3.53%  '
Set-TextCell 'D11' '1.017'
Set-TextCell 'E11' '  +5.00%  '
Set-TextCell 'D12' '22.38'
Set-TextCell 'E12' '  +6.55%  '
Set-TextCell 'B13' 'WrappedEther'
Set-TextCell 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell 'D13' '1.894.32'
Set-TextCell 'E13' '  +4.26%  '
Set-TextCell 'B14' 'Polkadot'
Set-TextCell 'C14' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextCell 'D14' '6.084'
Set-TextCell 'E14' '  +3.47%  '
Set-TextCell 'D15' '7.341'
Set-TextCell 'E15' '  +4.22%  '
Set-TextCell 'D16' '91.22'
Set-TextCell 'D17' '1.000'
Set-TextCell 'E17' '  -0.64%  '
Set-TextCell 'D18' '0.00001052'
Set-TextCell 'E18' '  +2.50%  '
Set-TextCell 'D19' '0.06612'
Set-TextCell 'E19' '  -0.53%  '
Set-TextCell 'D20' '17.77'
Set-TextCell 'E20' '  +4.04%  '
Set-TextCell 'D21' '1.001'
Set-TextCell 'E21' '  -0.40%  '
Set-TextCell 'D22' '28.889.56'
Set-TextCell 'E22' '  +5.30%  '
Set-TextCell 'D23' '5.552'
Set-TextCell 'E23' '  +3.98%  '
Set-TextCell 'E24' '  +3.12%  '
Set-TextCell 'E25' '  -1.66%  '
Set-TextCell 'D26' '2.114.54'
Set-TextCell 'E26' '  +3.86%  '
Set-TextCell 'D27' '161.24'
Set-TextCell 'E27' '  +3.64%  '
Set-TextCell 'D28' '19.99'
Set-TextCell 'E28' '  +3.22%  '
Set-TextCell 'D29' '2.162'
Set-TextCell 'E29' '  +5.35%  '
Set-TextCell 'D30' '5.529'
Set-TextCell 'E30' '  +4.97%  '
Set-TextCell 'D31' '120.37'
Set-TextCell 'E31' '  +1.73%  '
Set-TextCell 'D32' '1.011'
Set-TextCell 'E32' '  +7.29%  '
Set-TextCell 'D33' '0.09569'
Set-TextCell 'E33' '  +2.81%  '
Set-TextCell 'D34' '3.655'
Set-TextCell 'D35' '1.404'
Set-TextCell 'E35' '  +6.62%  '
Set-TextCell 'E36' '  +2.96%  '
Set-TextCell 'D37' '0.06198'
Set-TextCell 'E37' '  +4.56%  '
Set-TextCell 'D38' '0.02285'
Set-TextCell 'E38' '  +4.69%  '
Set-TextCell 'D39' '8.654'
Set-TextCell 'E39' '  +7.73%  '
Set-TextCell 'E40' '  +3.60%  '
Set-TextCell 'D41' '0.6011'
Set-TextCell 'E41' '  +4.54%  '
Set-TextCell 'D42' '0.1900'
Set-TextCell 'E42' '  +4.13%  '
Set-TextCell 'D43' '0.9997'
Set-TextCell 'E43' '  -0.52%  '
Set-TextCell 'D44' '10.48'
Set-TextCell 'E44' '  +5.08%  '
Set-TextCell 'D45' '1.268'
Set-TextCell 'E45' '  -0.01%  '
Set-TextCell 'D46' '0.5624'
Set-TextCell 'E46' '  +3.37%  '
Set-TextCell 'D47' '12.32'
Set-TextCell 'E47' '  +2.84%  '
Set-TextCell 'E48' '  +6.18%  '
Set-TextCell 'D49' '0.07255'
Set-TextCell 'E49' '  +9.81%  '
Set-TextCell 'D50' '2.122'
Set-TextCell 'E50' '  +16.52%  '
Set-TextCell 'D51' '112.95'
Set-TextCell 'E51' '  +1.96%  '
